$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -17.22788852866641
$ws.Cells.Item(2, 3).Value = -0.36098393029867
$ws.Cells.Item(2, 4).Value = -17.22788852866641
$ws.Cells.Item(2, 5).Value = -17.22788852866641
$ws.Cells.Item(2, 6).Value = -17.22788852866641
$ws.Cells.Item(2, 7).Value = -17.22788852866641
$ws.Cells.Item(2, 8).Value = -17.22788852866641
$ws.Cells.Item(2, 9).Value = -17.22788852866641
$ws.Cells.Item(2, 10).Value = -17.22788852866641
$ws.Cells.Item(2, 11).Value = -17.22788852866641
$ws.Cells.Item(3, 2).Value = -17.22788852866641
$ws.Cells.Item(3, 3).Value = -17.22788852866641
$ws.Cells.Item(3, 4).Value = -17.22788852866641
$ws.Cells.Item(3, 5).Value = -17.22788852866641
$ws.Cells.Item(3, 6).Value = -17.22788852866641
$ws.Cells.Item(3, 7).Value = -17.22788852866641
$ws.Cells.Item(3, 8).Value = -17.22788852866641
$ws.Cells.Item(3, 9).Value = 0.2923387960189974
$ws.Cells.Item(3, 10).Value = -17.22788852866641
$ws.Cells.Item(3, 11).Value = -17.22788852866641
$ws.Cells.Item(4, 2).Value = -17.22788852866641
$ws.Cells.Item(4, 3).Value = -0.4513912385021704
$ws.Cells.Item(4, 4).Value = 0.1191877818746619
$ws.Cells.Item(4, 5).Value = -17.22788852866641
$ws.Cells.Item(4, 6).Value = 3.963344827135121
$ws.Cells.Item(4, 7).Value = -17.22788852866641
$ws.Cells.Item(4, 8).Value = -17.22788852866641
$ws.Cells.Item(4, 9).Value = -17.22788852866641
$ws.Cells.Item(4, 10).Value = 2.516171694269409
$ws.Cells.Item(4, 11).Value = -17.22788852866641
$ws.Cells.Item(5, 2).Value = -17.22788852866641
$ws.Cells.Item(5, 3).Value = 0.3297621391670751
$ws.Cells.Item(5, 4).Value = -17.22788852866641
$ws.Cells.Item(5, 5).Value = -17.22788852866641
$ws.Cells.Item(5, 6).Value = -17.22788852866641
$ws.Cells.Item(5, 7).Value = 3.649588086775128
$ws.Cells.Item(5, 8).Value = -17.22788852866641
$ws.Cells.Item(5, 9).Value = -17.22788852866641
$ws.Cells.Item(5, 10).Value = -17.22788852866641
$ws.Cells.Item(5, 11).Value = -17.22788852866641
$ws.Cells.Item(6, 2).Value = -17.22788852866641
$ws.Cells.Item(6, 3).Value = -17.22788852866641
$ws.Cells.Item(6, 4).Value = -17.22788852866641
$ws.Cells.Item(6, 5).Value = -17.22788852866641
$ws.Cells.Item(6, 6).Value = -17.22788852866641
$ws.Cells.Item(6, 7).Value = -17.22788852866641
$ws.Cells.Item(6, 8).Value = -17.22788852866641
$ws.Cells.Item(6, 9).Value = -17.22788852866641
$ws.Cells.Item(6, 10).Value = -17.22788852866641
$ws.Cells.Item(6, 11).Value = -17.22788852866641
$ws.Cells.Item(7, 2).Value = 3.097602495966713
$ws.Cells.Item(7, 3).Value = -17.22788852866641
$ws.Cells.Item(7, 4).Value = -17.22788852866641
$ws.Cells.Item(7, 5).Value = -17.22788852866641
$ws.Cells.Item(7, 6).Value = -17.22788852866641
$ws.Cells.Item(7, 7).Value = -17.22788852866641
$ws.Cells.Item(7, 8).Value = -17.22788852866641
$ws.Cells.Item(7, 9).Value = -17.22788852866641
$ws.Cells.Item(7, 10).Value = -17.22788852866641
$ws.Cells.Item(7, 11).Value = -17.22788852866641
$ws.Cells.Item(8, 2).Value = -17.22788852866641
$ws.Cells.Item(8, 3).Value = -17.22788852866641
$ws.Cells.Item(8, 4).Value = -17.22788852866641
$ws.Cells.Item(8, 5).Value = 1.230393913375454
$ws.Cells.Item(8, 6).Value = -17.22788852866641
$ws.Cells.Item(8, 7).Value = -17.22788852866641
$ws.Cells.Item(8, 8).Value = -17.22788852866641
$ws.Cells.Item(8, 9).Value = -17.22788852866641
$ws.Cells.Item(8, 10).Value = -17.22788852866641
$ws.Cells.Item(8, 11).Value = -17.22788852866641
$ws.Cells.Item(9, 2).Value = 3.516006569930889
$ws.Cells.Item(9, 3).Value = -17.22788852866641
$ws.Cells.Item(9, 4).Value = -17.22788852866641
$ws.Cells.Item(9, 5).Value = -17.22788852866641
$ws.Cells.Item(9, 6).Value = -17.22788852866641
$ws.Cells.Item(9, 7).Value = -17.22788852866641
$ws.Cells.Item(9, 8).Value = -17.22788852866641
$ws.Cells.Item(9, 9).Value = -17.22788852866641
$ws.Cells.Item(9, 10).Value = -17.22788852866641
$ws.Cells.Item(9, 11).Value = -17.22788852866641
$ws.Cells.Item(10, 2).Value = -17.22788852866641
$ws.Cells.Item(10, 3).Value = -17.22788852866641
$ws.Cells.Item(10, 4).Value = -17.22788852866641
$ws.Cells.Item(10, 5).Value = -17.22788852866641
$ws.Cells.Item(10, 6).Value = -17.22788852866641
$ws.Cells.Item(10, 7).Value = -17.22788852866641
$ws.Cells.Item(10, 8).Value = -17.22788852866641
$ws.Cells.Item(10, 9).Value = 0.9297001910298794
$ws.Cells.Item(10, 10).Value = -17.22788852866641
$ws.Cells.Item(10, 11).Value = 1.818769231911305
$ws.Cells.Item(11, 2).Value = -17.22788852866641
$ws.Cells.Item(11, 3).Value = -17.22788852866641
$ws.Cells.Item(11, 4).Value = -17.22788852866641
$ws.Cells.Item(11, 5).Value = 2.323123102198249
$ws.Cells.Item(11, 6).Value = -17.22788852866641
$ws.Cells.Item(11, 7).Value = 1.077037595870834
$ws.Cells.Item(11, 8).Value = -17.22788852866641
$ws.Cells.Item(11, 9).Value = -17.22788852866641
$ws.Cells.Item(11, 10).Value = -17.22788852866641
$ws.Cells.Item(11, 11).Value = 1.601651504674222
$ws.Cells.Item(12, 2).Value = -17.22788852866641
$ws.Cells.Item(12, 3).Value = -17.22788852866641
$ws.Cells.Item(12, 4).Value = -17.22788852866641
$ws.Cells.Item(12, 5).Value = -17.22788852866641
$ws.Cells.Item(12, 6).Value = -17.22788852866641
$ws.Cells.Item(12, 7).Value = -17.22788852866641
$ws.Cells.Item(12, 8).Value = -17.22788852866641
$ws.Cells.Item(12, 9).Value = -17.22788852866641
$ws.Cells.Item(12, 10).Value = -17.22788852866641
$ws.Cells.Item(12, 11).Value = -17.22788852866641
$ws.Cells.Item(13, 2).Value = -17.22788852866641
$ws.Cells.Item(13, 3).Value = -17.22788852866641
$ws.Cells.Item(13, 4).Value = -17.22788852866641
$ws.Cells.Item(13, 5).Value = 2.209207202823092
$ws.Cells.Item(13, 6).Value = -17.22788852866641
$ws.Cells.Item(13, 7).Value = -17.22788852866641
$ws.Cells.Item(13, 8).Value = -17.22788852866641
$ws.Cells.Item(13, 9).Value = -17.22788852866641
$ws.Cells.Item(13, 10).Value = 1.39081599807974
$ws.Cells.Item(13, 11).Value = 2.866432758281237
$ws.Cells.Item(14, 2).Value = -17.22788852866641
$ws.Cells.Item(14, 3).Value = -17.22788852866641
$ws.Cells.Item(14, 4).Value = 0.8073240406011979
$ws.Cells.Item(14, 5).Value = -17.22788852866641
$ws.Cells.Item(14, 6).Value = -17.22788852866641
$ws.Cells.Item(14, 7).Value = -17.22788852866641
$ws.Cells.Item(14, 8).Value = -17.22788852866641
$ws.Cells.Item(14, 9).Value = -17.22788852866641
$ws.Cells.Item(14, 10).Value = -17.22788852866641
$ws.Cells.Item(14, 11).Value = 2.140867289776827
$ws.Cells.Item(15, 2).Value = -17.22788852866641
$ws.Cells.Item(15, 3).Value = -17.22788852866641
$ws.Cells.Item(15, 4).Value = -0.05147185959816603
$ws.Cells.Item(15, 5).Value = -17.22788852866641
$ws.Cells.Item(15, 6).Value = -17.22788852866641
$ws.Cells.Item(15, 7).Value = -17.22788852866641
$ws.Cells.Item(15, 8).Value = -17.22788852866641
$ws.Cells.Item(15, 9).Value = -17.22788852866641
$ws.Cells.Item(15, 10).Value = -17.22788852866641
$ws.Cells.Item(15, 11).Value = -17.22788852866641
$ws.Cells.Item(16, 2).Value = -17.22788852866641
$ws.Cells.Item(16, 3).Value = -17.22788852866641
$ws.Cells.Item(16, 4).Value = -17.22788852866641
$ws.Cells.Item(16, 5).Value = -17.22788852866641
$ws.Cells.Item(16, 6).Value = -17.22788852866641
$ws.Cells.Item(16, 7).Value = -17.22788852866641
$ws.Cells.Item(16, 8).Value = -17.22788852866641
$ws.Cells.Item(16, 9).Value = -17.22788852866641
$ws.Cells.Item(16, 10).Value = 1.347447588526066
$ws.Cells.Item(16, 11).Value = -17.22788852866641
$ws.Cells.Item(17, 2).Value = -17.22788852866641
$ws.Cells.Item(17, 3).Value = 0.9320665504083744
$ws.Cells.Item(17, 4).Value = -0.01113319902740253
$ws.Cells.Item(17, 5).Value = -17.22788852866641
$ws.Cells.Item(17, 6).Value = -17.22788852866641
$ws.Cells.Item(17, 7).Value = -17.22788852866641
$ws.Cells.Item(17, 8).Value = -17.22788852866641
$ws.Cells.Item(17, 9).Value = 1.064509500811042
$ws.Cells.Item(17, 10).Value = 2.215048894402827
$ws.Cells.Item(17, 11).Value = -17.22788852866641
$ws.Cells.Item(18, 2).Value = -17.22788852866641
$ws.Cells.Item(18, 3).Value = -17.22788852866641
$ws.Cells.Item(18, 4).Value = -17.22788852866641
$ws.Cells.Item(18, 5).Value = -17.22788852866641
$ws.Cells.Item(18, 6).Value = -17.22788852866641
$ws.Cells.Item(18, 7).Value = -17.22788852866641
$ws.Cells.Item(18, 8).Value = -17.22788852866641
$ws.Cells.Item(18, 9).Value = -0.5069362450747648
$ws.Cells.Item(18, 10).Value = 2.160088802758856
$ws.Cells.Item(18, 11).Value = -17.22788852866641
$ws.Cells.Item(19, 2).Value = -17.22788852866641
$ws.Cells.Item(19, 3).Value = -17.22788852866641
$ws.Cells.Item(19, 4).Value = 3.065314790224949
$ws.Cells.Item(19, 5).Value = -17.22788852866641
$ws.Cells.Item(19, 6).Value = -17.22788852866641
$ws.Cells.Item(19, 7).Value = -17.22788852866641
$ws.Cells.Item(19, 8).Value = 4.321919166186377
$ws.Cells.Item(19, 9).Value = 1.289338714189126
$ws.Cells.Item(19, 10).Value = -17.22788852866641
$ws.Cells.Item(19, 11).Value = -17.22788852866641
$ws.Cells.Item(20, 2).Value = -17.22788852866641
$ws.Cells.Item(20, 3).Value = 3.182311251826417
$ws.Cells.Item(20, 4).Value = 2.773163621399769
$ws.Cells.Item(20, 5).Value = -17.22788852866641
$ws.Cells.Item(20, 6).Value = 2.137923409335945
$ws.Cells.Item(20, 7).Value = -17.22788852866641
$ws.Cells.Item(20, 8).Value = -17.22788852866641
$ws.Cells.Item(20, 9).Value = 3.53991201864155
$ws.Cells.Item(20, 10).Value = -17.22788852866641
$ws.Cells.Item(20, 11).Value = 0.7943763239356866
$ws.Cells.Item(21, 2).Value = -17.22788852866641
$ws.Cells.Item(21, 3).Value = 2.643284265387411
$ws.Cells.Item(21, 4).Value = -17.22788852866641
$ws.Cells.Item(21, 5).Value = 3.00454293873457
$ws.Cells.Item(21, 6).Value = -17.22788852866641
$ws.Cells.Item(21, 7).Value = 2.416955301387179
$ws.Cells.Item(21, 8).Value = -17.22788852866641
$ws.Cells.Item(21, 9).Value = -17.22788852866641
$ws.Cells.Item(21, 10).Value = -17.22788852866641
$ws.Cells.Item(21, 11).Value = -17.22788852866641
